$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("December")

$ws.Range("B2").Value = 1324
$ws.Range("C2").Value = 1102
$ws.Range("D2").Value = 222
$ws.Range("E2").Value = "We borrowerd more than we lent"
$ws.Range("G2").Value = "1.20 : 1"

$ws.Range("B3").Value = 510
$ws.Range("C3").Value = 369
$ws.Range("D3").Value = 141
$ws.Range("E3").Value = "We borrowerd more than we lent"
$ws.Range("G3").Value = "1.38 : 1"

$ws.Range("B4").Value = 975
$ws.Range("C4").Value = 1098
$ws.Range("D4").Value = -123
$ws.Range("F4").Value = "We lent more than we borrowed"
$ws.Range("G4").Value = "0.89 : 1"

$ws.Range("B5").Value = 29
$ws.Range("C5").Value = 104
$ws.Range("D5").Value = -75
$ws.Range("F5").Value = "We lent more than we borrowed"
$ws.Range("G5").Value = "0.28 : 1"

$ws.Range("B6").Value = 836
$ws.Range("C6").Value = 1375
$ws.Range("D6").Value = -539
$ws.Range("F6").Value = "We lent more than we borrowed"
$ws.Range("G6").Value = "0.61 : 1"

$ws.Range("B7").Value = 154
$ws.Range("C7").Value = 158
$ws.Range("D7").Value = -4
$ws.Range("F7").Value = "We lent more than we borrowed"
$ws.Range("G7").Value = "0.97 : 1"

$ws.Range("B8").Value = 140
$ws.Range("C8").Value = 143
$ws.Range("D8").Value = -3
$ws.Range("F8").Value = "We lent more than we borrowed"
$ws.Range("G8").Value = "0.98 : 1"

$ws.Range("B9").Value = 32
$ws.Range("C9").Value = 49
$ws.Range("D9").Value = -17
$ws.Range("F9").Value = "We lent more than we borrowed"
$ws.Range("G9").Value = "0.65 : 1"

$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 26
$ws.Range("D10").Value = -23
$ws.Range("F10").Value = "We lent more than we borrowed"
$ws.Range("G10").Value = "0.12 : 1"

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 19
$ws.Range("C12").Value = 17
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = "We borrowerd more than we lent"
$ws.Range("G12").Value = "1.12 : 1"

$ws.Range("B13").Value = 91
$ws.Range("C13").Value = 54
$ws.Range("D13").Value = 37
$ws.Range("E13").Value = "We borrowerd more than we lent"
$ws.Range("G13").Value = "1.69 : 1"

$ws.Range("B14").Value = 144
$ws.Range("C14").Value = 159
$ws.Range("D14").Value = -15
$ws.Range("F14").Value = "We lent more than we borrowed"
$ws.Range("G14").Value = "0.91 : 1"

$ws.Range("B15").Value = 43
$ws.Range("C15").Value = 133
$ws.Range("D15").Value = -90
$ws.Range("F15").Value = "We lent more than we borrowed"
$ws.Range("G15").Value = "0.32 : 1"

$ws.Range("B16").Value = 64
$ws.Range("C16").Value = 147
$ws.Range("D16").Value = -83
$ws.Range("F16").Value = "We lent more than we borrowed"
$ws.Range("G16").Value = "0.44 : 1"

$ws.Range("B17").Value = 499
$ws.Range("C17").Value = 387
$ws.Range("D17").Value = 112
$ws.Range("E17").Value = "We borrowerd more than we lent"
$ws.Range("G17").Value = "1.29 : 1"

$ws.Range("B18").Value = 64
$ws.Range("C18").Value = 61
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = "We borrowerd more than we lent"
$ws.Range("G18").Value = "1.05 : 1"

$ws.Range("B19").Value = 427
$ws.Range("C19").Value = 344
$ws.Range("D19").Value = 83
$ws.Range("E19").Value = "We borrowerd more than we lent"
$ws.Range("G19").Value = "1.24 : 1"

$ws.Range("B20").Value = 20
$ws.Range("C20").Value = 38
$ws.Range("D20").Value = -18
$ws.Range("F20").Value = "We lent more than we borrowed"
$ws.Range("G20").Value = "0.53 : 1"

$ws.Range("B21").Value = 477
$ws.Range("C21").Value = 260
$ws.Range("D21").Value = 217
$ws.Range("E21").Value = "We borrowerd more than we lent"
$ws.Range("G21").Value = "1.83 : 1"

$ws.Range("B22").Value = 45
$ws.Range("C22").Value = 153
$ws.Range("D22").Value = -108
$ws.Range("F22").Value = "We lent more than we borrowed"
$ws.Range("G22").Value = "0.29 : 1"

$ws.Range("B23").Value = 467
$ws.Range("C23").Value = 264
$ws.Range("D23").Value = 203
$ws.Range("E23").Value = "We borrowerd more than we lent"
$ws.Range("G23").Value = "1.77 : 1"

$ws.Range("B24").Value = 1420
$ws.Range("C24").Value = 1060
$ws.Range("D24").Value = 360
$ws.Range("E24").Value = "We borrowerd more than we lent"
$ws.Range("G24").Value = "1.34 : 1"

$ws.Range("B25").Value = 177
$ws.Range("C25").Value = 397
$ws.Range("D25").Value = -220
$ws.Range("F25").Value = "We lent more than we borrowed"
$ws.Range("G25").Value = "0.45 : 1"

$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0

$ws.Range("B27").Value = 227
$ws.Range("C27").Value = 147
$ws.Range("D27").Value = 80
$ws.Range("E27").Value = "We borrowerd more than we lent"
$ws.Range("G27").Value = "1.54 : 1"

$ws.Range("B28").Value = 73
$ws.Range("C28").Value = 87
$ws.Range("D28").Value = -14
$ws.Range("F28").Value = "We lent more than we borrowed"
$ws.Range("G28").Value = "0.84 : 1"

$ws.Range("B29").Value = 539
$ws.Range("C29").Value = 441
$ws.Range("D29").Value = 98
$ws.Range("E29").Value = "We borrowerd more than we lent"
$ws.Range("G29").Value = "1.22 : 1"

$ws.Range("B30").Value = 55
$ws.Range("C30").Value = 28
$ws.Range("D30").Value = 27
$ws.Range("E30").Value = "We borrowerd more than we lent"
$ws.Range("G30").Value = "1.96 : 1"

$ws.Range("B31").Value = 47
$ws.Range("C31").Value = 171
$ws.Range("D31").Value = -124
$ws.Range("F31").Value = "We lent more than we borrowed"
$ws.Range("G31").Value = "0.27 : 1"

$ws.Range("B32").Value = 354
$ws.Range("C32").Value = 503
$ws.Range("D32").Value = -149
$ws.Range("F32").Value = "We lent more than we borrowed"
$ws.Range("G32").Value = "0.70 : 1"

$ws.Range("B33").Value = 322
$ws.Range("C33").Value = 427
$ws.Range("D33").Value = -105
$ws.Range("F33").Value = "We lent more than we borrowed"
$ws.Range("G33").Value = "0.75 : 1"

$ws.Range("B34").Value = 178
$ws.Range("C34").Value = 117
$ws.Range("D34").Value = 61
$ws.Range("E34").Value = "We borrowerd more than we lent"
$ws.Range("G34").Value = "1.52 : 1"

$ws.Range("B35").Value = 732
$ws.Range("C35").Value = 997
$ws.Range("D35").Value = -265
$ws.Range("F35").Value = "We lent more than we borrowed"
$ws.Range("G35").Value = "0.73 : 1"

$ws.Range("B36").Value = 123
$ws.Range("C36").Value = 439
$ws.Range("D36").Value = -316
$ws.Range("F36").Value = "We lent more than we borrowed"
$ws.Range("G36").Value = "0.28 : 1"

$ws.Range("B37").Value = 484
$ws.Range("C37").Value = 308
$ws.Range("D37").Value = 176
$ws.Range("E37").Value = "We borrowerd more than we lent"
$ws.Range("G37").Value = "1.57 : 1"

$ws.Range("B38").Value = 31
$ws.Range("C38").Value = 129
$ws.Range("D38").Value = -98
$ws.Range("F38").Value = "We lent more than we borrowed"
$ws.Range("G38").Value = "0.24 : 1"

$ws.Range("B39").Value = 21
$ws.Range("C39").Value = 80
$ws.Range("D39").Value = -59
$ws.Range("F39").Value = "We lent more than we borrowed"
$ws.Range("G39").Value = "0.26 : 1"

$ws.Range("B40").Value = 88
$ws.Range("C40").Value = 146
$ws.Range("D40").Value = -58
$ws.Range("F40").Value = "We lent more than we borrowed"
$ws.Range("G40").Value = "0.60 : 1"

$ws.Range("B41").Value = 3
$ws.Range("C41").Value = 26
$ws.Range("D41").Value = -23
$ws.Range("F41").Value = "We lent more than we borrowed"
$ws.Range("G41").Value = "0.12 : 1"

$ws.Range("B42").Value = 6
$ws.Range("C42").Value = 23
$ws.Range("D42").Value = -17
$ws.Range("F42").Value = "We lent more than we borrowed"
$ws.Range("G42").Value = "0.26 : 1"

$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0

$ws.Range("B44").Value = 74
$ws.Range("C44").Value = 75
$ws.Range("D44").Value = -1
$ws.Range("F44").Value = "We lent more than we borrowed"
$ws.Range("G44").Value = "0.99 : 1"

$ws.Range("B45").Value = 63
$ws.Range("C45").Value = 171
$ws.Range("D45").Value = -108
$ws.Range("F45").Value = "We lent more than we borrowed"
$ws.Range("G45").Value = "0.37 : 1"

$ws.Range("B46").Value = 399
$ws.Range("C46").Value = 389
$ws.Range("D46").Value = 10
$ws.Range("E46").Value = "We borrowerd more than we lent"
$ws.Range("G46").Value = "1.03 : 1"

$ws.Range("B47").Value = 792
$ws.Range("C47").Value = 506
$ws.Range("D47").Value = 286
$ws.Range("E47").Value = "We borrowerd more than we lent"
$ws.Range("G47").Value = "1.57 : 1"

$ws.Range("B48").Value = 155
$ws.Range("C48").Value = 537
$ws.Range("D48").Value = -382
$ws.Range("F48").Value = "We lent more than we borrowed"
$ws.Range("G48").Value = "0.29 : 1"

$ws.Range("B49").Value = 650
$ws.Range("C49").Value = 184
$ws.Range("D49").Value = 466
$ws.Range("E49").Value = "We borrowerd more than we lent"
$ws.Range("G49").Value = "3.53 : 1"

$ws.Range("B50").Value = 777
$ws.Range("C50").Value = 471
$ws.Range("D50").Value = 306
$ws.Range("E50").Value = "We borrowerd more than we lent"
$ws.Range("G50").Value = "1.65 : 1"

$ws.Range("B51").Value = 161
$ws.Range("C51").Value = 171
$ws.Range("D51").Value = -10
$ws.Range("F51").Value = "We lent more than we borrowed"
$ws.Range("G51").Value = "0.94 : 1"

$ws.Range("B52").Value = 209
$ws.Range("C52").Value = 226
$ws.Range("D52").Value = -17
$ws.Range("F52").Value = "We lent more than we borrowed"
$ws.Range("G52").Value = "0.92 : 1"

$ws.Range("B53").Value = 132
$ws.Range("C53").Value = 165
$ws.Range("D53").Value = -33
$ws.Range("F53").Value = "We lent more than we borrowed"
$ws.Range("G53").Value = "0.80 : 1"

$ws.Range("B54").Value = 30
$ws.Range("C54").Value = 197
$ws.Range("D54").Value = -167
$ws.Range("F54").Value = "We lent more than we borrowed"
$ws.Range("G54").Value = "0.15 : 1"

$ws.Range("B55").Value = 519
$ws.Range("C55").Value = 145
$ws.Range("D55").Value = 374
$ws.Range("E55").Value = "We borrowerd more than we lent"
$ws.Range("G55").Value = "3.58 : 1"

$ws.Activate()
$ws.Range("B2").Select()